$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values updated (recalculated means)
$ws.Range("B2").Value = -0.7154229613570755
$ws.Range("C2").Value = 5.2025224677145037
$ws.Range("D2").Value = 6.8414277798401288
$ws.Range("E2").Value = 11.54005608882097

# Row 3 data values updated; C3 cleared, D3 gains a new value
$ws.Range("B3").Value = -6.2126054779321409
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 14.078340384993421
$ws.Range("E3").Value = 5.8348200070450318

# Update the active selection to match the new range used (B1:E3)
$ws.Range("B1:E3").Select()
